$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D5 content: was [{"money":"100"}] -> now [{"hero":{"id":"2"}}]
$ws.Range("D5").Value = '[{"hero":{"id":"2"}}]'

# Update the active selection to D5 (as seen in sheetView selection change)
$ws.Range("D5").Select()
